$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 96, shifting existing rows 96:194 down to 97:195
$ws.Rows("96:96").Insert()

# Populate the newly inserted row 96 with data
# (a record for Vega Modelo de Temuco / Pomelo / Start Ruby / Primera)
$ws.Range("A96").Value = 10
$ws.Range("B96").Value = "Vega Modelo de Temuco"
$ws.Range("C96").Value = "La Araucanía"
$ws.Range("D96").Value = 44601
$ws.Range("E96").Value = 9
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100102
$ws.Range("H96").Value = "Cítricos"
$ws.Range("I96").Value = 100102006
$ws.Range("J96").Value = "Pomelo"
$ws.Range("K96").Value = "Start Ruby"
$ws.Range("L96").Value = "Primera"
$ws.Range("M96").Value = 50
$ws.Range("N96").Value = 13000
$ws.Range("O96").Value = 13000
$ws.Range("P96").Value = 13000
$ws.Range("Q96").Value = "$/bandeja 15 kilos granel"
$ws.Range("R96").Value = "Región de O'Higgins"
$ws.Range("S96").Value = 867
$ws.Range("T96").Value = 15
